# Add a new trade row (row 6) to the GILD trade log sheet, matching the
# columns/style already used by the existing trade rows (3-5), and widen
# column E (BuyPrice) to fit the new, wider value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new trade -------------------------------------------------
# Copy formatting from the row above (row 5) so the new row picks up the
# same number/boolean styles already used by the other trade rows, then
# overwrite with the new values.
$ws.Range("A5:I5").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)

$ws.Range("A6").Value = 42649.654166666667
$ws.Range("B6").Value = $false
$ws.Range("C6").Value = 10010.959999999999
$ws.Range("D6").Value = 10015.469999999999
$ws.Range("E6").Value = 77.349997999999999
$ws.Range("F6").Value = 77.42
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 0.09
$ws.Range("I6").Value = $false

# --- Column E width -----------------------------------------------------
# The new BuyPrice value needs a slightly wider column to continue
# best-fitting its contents.
$ws.Columns("E").ColumnWidth = 9
